# Autogenerated on Wed Apr 01 2015 00:15:40 GMT+0000 (Coordinated Universal Time)
# Update MSME Country Indicators - Slovak Republic Summary figures with
# more precise (2-decimal) percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Leading apostrophes keep these percentage-like figures stored as text,
# matching how the original workbook stores them (shared strings, not
# numeric cells).

# Row 33: Enterprises density (per 1000 people)
$ws.Range("B33").Value = "'73.78"
$ws.Range("C33").Value = "'2.87"
$ws.Range("D33").Value = "'76.65"

# Row 34: Employment (% of total)
$ws.Range("B34").Value = "'38.48"
$ws.Range("D34").Value = "'71.68"

# Row 36: Enterprises (% of total)
$ws.Range("B36").Value = "'96.15"
$ws.Range("C36").Value = "'3.73"
$ws.Range("D36").Value = "'99.88"

# Row 40: Value added to the economy (% of total)
$ws.Range("B40").Value = "'26.44"
$ws.Range("C40").Value = "'37.63"
$ws.Range("D40").Value = "'64.07"
